# Auto-update draw results: append the newest Pick 3 draw as a new row
# at the bottom of the results table on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the next blank row right after the current table (row 29 -> 30).
$lastRow = $ws.UsedRange.Rows.Count
$newRow  = $lastRow + 1

$date        = "2025-10-16"
$game        = "Pick 3"
$phase       = "251016"
$result      = "6-0-1"
$insertedAt  = "2025-10-16T21:38:22.187+04:00"

# Date- and digit-only values must be forced to Text so Excel doesn't
# silently reinterpret them as a date serial / number (the source data
# is plain text in every column).
$ws.Range("A" + $newRow).NumberFormat = "@"
$ws.Range("C" + $newRow).NumberFormat = "@"

$ws.Cells.Item($newRow, 1).Value = $date
$ws.Cells.Item($newRow, 2).Value = $game
$ws.Cells.Item($newRow, 3).Value = $phase
$ws.Cells.Item($newRow, 4).Value = $result
$ws.Cells.Item($newRow, 5).Value = $insertedAt

# Keep the "number stored as text" warning suppressed across the whole
# (now one row taller) table, same as the rest of the sheet.
$fullRange = $ws.Range("A1:E" + $newRow)
$fullRange.Errors.Item(3).Ignore = $true
